$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.284.17"
Set-TextValue "E2" "  -2.89%  "
Set-TextValue "D3" "1.854.65"
Set-TextValue "E3" "  -3.67%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.18%  "
Set-TextValue "D5" "324.02"
Set-TextValue "E5" "  -1.50%  "
Set-TextValue "E6" "  -0.17%  "
Set-TextValue "D7" "0.4557"
Set-TextValue "E7" "  -3.56%  "
Set-TextValue "D8" "0.3870"
Set-TextValue "E8" "  -4.84%  "
Set-TextValue "D9" "48.42"
Set-TextValue "E9" "  -8.72%  "
Set-TextValue "D10" "0.07930"
Set-TextValue "E10" "  -6.04%  "
Set-TextValue "D11" "1.015"
Set-TextValue "E11" "  -3.07%  "
Set-TextValue "E12" "  -4.02%  "
Set-TextValue "D13" "1.857.07"
Set-TextValue "E13" "  -3.31%  "
Set-TextValue "D14" "5.899"
Set-TextValue "E14" "  -3.30%  "
Set-TextValue "D15" "7.147"
Set-TextValue "E15" "  -4.87%  "
Set-TextValue "E16" "  -0.33%  "
Set-TextValue "D17" "85.85"
Set-TextValue "E17" "  -5.48%  "
Set-TextValue "E18" "  -0.55%  "
Set-TextValue "D19" "0.00001023"
Set-TextValue "E19" "  -4.08%  "
Set-TextValue "D20" "17.07"
Set-TextValue "E20" "  -5.56%  "
Set-TextValue "D21" "1.001"
Set-TextValue "E21" "  -0.13%  "
Set-TextValue "D22" "5.502"
Set-TextValue "E22" "  -4.48%  "
Set-TextValue "D23" "27.285.67"
Set-TextValue "E23" "  -2.89%  "
Set-TextValue "D24" "10.89"
Set-TextValue "E24" "  -4.69%  "
Set-TextValue "D25" "2.288"
Set-TextValue "E25" "  +0.33%  "
Set-TextValue "D26" "2.069.69"
Set-TextValue "E26" "  -4.51%  "
Set-TextValue "D27" "153.76"
Set-TextValue "E27" "  -0.45%  "
Set-TextValue "D28" "19.89"
Set-TextValue "E28" "  -1.19%  "
Set-TextValue "E29" "  -4.52%  "
Set-TextValue "D30" "5.461"
Set-TextValue "E30" "  -4.54%  "
Set-TextValue "D31" "121.19"
Set-TextValue "E31" "  -2.12%  "
Set-TextValue "D32" "0.09322"
Set-TextValue "E32" "  -2.90%  "
Set-TextValue "D33" "0.9363"
Set-TextValue "E33" "  -3.94%  "
Set-TextValue "D34" "1.455"
Set-TextValue "E34" "  +0.87%  "
Set-TextValue "D35" "3.584"
Set-TextValue "E35" "  -1.56%  "
Set-TextValue "D36" "5.269"
Set-TextValue "E36" "  -5.22%  "
Set-TextValue "D37" "0.02224"
Set-TextValue "E37" "  -4.05%  "
Set-TextValue "D38" "0.05999"
Set-TextValue "E38" "  -2.72%  "
Set-TextValue "D39" "1.219"
Set-TextValue "E39" "  -1.51%  "
Set-TextValue "D40" "8.047"
Set-TextValue "E40" "  -10.99%  "
Set-TextValue "E41" "  -0.21%  "
Set-TextValue "D42" "0.5909"
Set-TextValue "E42" "  -4.33%  "
Set-TextValue "D43" "0.1888"
Set-TextValue "E43" "  -1.00%  "
Set-TextValue "D44" "10.14"
Set-TextValue "E44" "  -8.40%  "
Set-TextValue "D45" "1.284"
Set-TextValue "E45" "  -0.51%  "
Set-TextValue "D46" "0.5604"
Set-TextValue "E46" "  -5.01%  "
Set-TextValue "D47" "12.00"
Set-TextValue "E47" "  -6.53%  "
Set-TextValue "D48" "3.372"
Set-TextValue "E48" "  -2.78%  "
Set-TextValue "D49" "1.916"
Set-TextValue "E49" "  -5.94%  "
Set-TextValue "D50" "0.06737"
Set-TextValue "E50" "  -1.20%  "
Set-TextValue "D51" "108.64"
Set-TextValue "E51" "  -1.28%  "
